# Update cryptocurrency price/volume symbol list (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the literal text into the cell (preserving the inline/shared-string
    # "t=str" type for numeric-looking text like prices and percentages) and then
    # restore the default "Normal" cell style so NumberFormat="@" does not leave
    # a stray Text-format style applied to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "277.73"
Set-TextValue $ws.Range("E2") "1.72%"

Set-TextValue $ws.Range("D3") "27.22"
Set-TextValue $ws.Range("E3") "1.62%"

Set-TextValue $ws.Range("D4") "4.870"
Set-TextValue $ws.Range("E4") "-0.84%"

Set-TextValue $ws.Range("D5") "0.06394"
Set-TextValue $ws.Range("E5") "1.28%"

Set-TextValue $ws.Range("D6") "6.982"
Set-TextValue $ws.Range("E6") "1.17%"

Set-TextValue $ws.Range("D7") "1.250"
Set-TextValue $ws.Range("E7") "-6.97%"

Set-TextValue $ws.Range("D8") "0.8820"
Set-TextValue $ws.Range("E8") "-0.21%"

Set-TextValue $ws.Range("D9") "0.1518"
Set-TextValue $ws.Range("E9") "3.09%"

Set-TextValue $ws.Range("D10") "0.05127"
Set-TextValue $ws.Range("E10") "0.64%"

Set-TextValue $ws.Range("D11") "0.07529"
Set-TextValue $ws.Range("E11") "1.87%"

Set-TextValue $ws.Range("D12") "0.02954"
Set-TextValue $ws.Range("E12") "-7.60%"

Set-TextValue $ws.Range("D13") "0.09015"
Set-TextValue $ws.Range("E13") "-0.40%"

Set-TextValue $ws.Range("D14") "0.001563"
Set-TextValue $ws.Range("E14") "-1.16%"

Set-TextValue $ws.Range("D15") "0.0006404"
Set-TextValue $ws.Range("E15") "1.13%"

Set-TextValue $ws.Range("D16") "0.005994"
Set-TextValue $ws.Range("E16") "-1.42%"

Set-TextValue $ws.Range("D17") "3.461"
Set-TextValue $ws.Range("E17") "-0.32%"

Set-TextValue $ws.Range("D18") "3.320"
Set-TextValue $ws.Range("E18") "-0.83%"

Set-TextValue $ws.Range("E19") "-0.51%"

Set-TextValue $ws.Range("E20") "1.35%"

Set-TextValue $ws.Range("E21") "-0.01%"

Set-TextValue $ws.Range("D22") "3.906"
Set-TextValue $ws.Range("E22") "-0.58%"

Set-TextValue $ws.Range("D23") "0.04421"
Set-TextValue $ws.Range("E23") "1.57%"

Set-TextValue $ws.Range("D24") "0.001173"
Set-TextValue $ws.Range("E24") "-0.68%"

Set-TextValue $ws.Range("D25") "0.003880"
Set-TextValue $ws.Range("E25") "6.62%"

Set-TextValue $ws.Range("E26") "-0.05%"

Set-TextValue $ws.Range("E27") "13.74%"

Set-TextValue $ws.Range("D40") "0.04160"
Set-TextValue $ws.Range("E40") "2.51%"

Set-TextValue $ws.Range("D41") "0.006853"
Set-TextValue $ws.Range("E41") "3.76%"

Set-TextValue $ws.Range("D42") "0.1179"
Set-TextValue $ws.Range("E42") "1.47%"

Set-TextValue $ws.Range("D43") "0.002032"
Set-TextValue $ws.Range("E43") "-5.62%"

Set-TextValue $ws.Range("D44") "0.01125"
Set-TextValue $ws.Range("E44") "-10.93%"

Set-TextValue $ws.Range("D45") "0.00005180"
Set-TextValue $ws.Range("E45") "-3.23%"

Set-TextValue $ws.Range("D46") "1.482"
Set-TextValue $ws.Range("E46") "-37.07%"

Set-TextValue $ws.Range("D47") "0.02024"
Set-TextValue $ws.Range("E47") "-4.80%"
